# Apply updated Leve profit values across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 162.75
$ws.Range("I6").Value = 108.666664
$ws.Range("J6").Value = 325
$ws.Range("K6").Value = 325.999992
$ws.Range("L6").Value = 975
$ws.Range("M6").Value = -213.999992
$ws.Range("N6").Value = -1199

$ws.Range("H29").Value = 933.3333
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 3500
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 10500
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -11062

$ws.Range("H55").Value = 181
$ws.Range("I55").Value = 135
$ws.Range("J55").Value = 257.66666
$ws.Range("K55").Value = 135
$ws.Range("L55").Value = 257.66666
$ws.Range("M55").Value = 79
$ws.Range("N55").Value = -685.66666

$ws.Range("H64").Value = 36502
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 36502
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 36502
$ws.Range("N64").Value = -36998

$ws.Range("H67").Value = 36502
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 36502
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 36502
$ws.Range("N67").Value = -38218

$ws.Range("H74").Value = 2500
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2500
$ws.Range("N74").Value = -4372

$ws.Range("H77").Value = 2500
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 12500
$ws.Range("N77").Value = -21860

$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -4008

$ws.Range("H132").Value = 845
$ws.Range("I132").Value = 815.4286
$ws.Range("J132").Value = 1052
$ws.Range("K132").Value = 2446.2858
$ws.Range("L132").Value = 3156
$ws.Range("M132").Value = 83.71420000000035
$ws.Range("N132").Value = -8216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14463.23
$ws.Range("I32").Value = 12802.2
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 12802.2
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -12515.2
$ws.Range("N32").Value = -20574

$ws.Range("H122").Value = 1748.6666
$ws.Range("I122").Value = 1748.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5245.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2795.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6145.75
$ws.Range("I86").Value = 2475
$ws.Range("J86").Value = 6879.9
$ws.Range("K86").Value = 2475
$ws.Range("L86").Value = 6879.9
$ws.Range("M86").Value = -1352
$ws.Range("N86").Value = -9125.9

$ws.Range("H89").Value = 6145.75
$ws.Range("I89").Value = 2475
$ws.Range("J89").Value = 6879.9
$ws.Range("K89").Value = 12375
$ws.Range("L89").Value = 34399.5
$ws.Range("M89").Value = -6759
$ws.Range("N89").Value = -45631.5

$ws.Range("H134").Value = 4035.3333
$ws.Range("I134").Value = 3856.6365
$ws.Range("J134").Value = 6001
$ws.Range("K134").Value = 11569.9095
$ws.Range("L134").Value = 18003
$ws.Range("M134").Value = -9034.9095
$ws.Range("N134").Value = -23073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 603.2857
$ws.Range("I19").Value = 603.2857
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 603.2857
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -433.2857
$ws.Range("N19").ClearContents()

$ws.Range("H24").Value = 603.2857
$ws.Range("I24").Value = 603.2857
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 603.2857
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -433.2857
$ws.Range("N24").ClearContents()

$ws.Range("H99").Value = 2863
$ws.Range("I99").Value = 2863
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2863
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1365

$ws.Range("H122").Value = 1971.1428
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 2998
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 8994
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -13894

$ws.Range("H126").Value = 2863
$ws.Range("I126").Value = 2863
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8589
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6119

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1930.6666
$ws.Range("I117").Value = 1132.25
$ws.Range("J117").Value = 2569.4
$ws.Range("K117").Value = 3396.75
$ws.Range("L117").Value = 7708.200000000001
$ws.Range("M117").Value = 45.25
$ws.Range("N117").Value = -14592.2

$ws.Range("H129").Value = 2080
$ws.Range("I129").Value = 2225
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 6675
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = -1675
$ws.Range("N129").Value = -14500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H107").Value = 296.83334
$ws.Range("I107").Value = 96.2
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 96.2
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 1823.8
$ws.Range("N107").Value = -5140

$ws.Range("H122").Value = 7414
$ws.Range("I122").Value = 3203.611
$ws.Range("J122").Value = 22571.4
$ws.Range("K122").Value = 9610.832999999999
$ws.Range("L122").Value = 67714.20000000001
$ws.Range("M122").Value = -7160.832999999999
$ws.Range("N122").Value = -72614.20000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 5340.6665
$ws.Range("I26").Value = 5340.6665
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 5340.6665
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -5045.6665
$ws.Range("N26").ClearContents()

$ws.Range("H122").Value = 6468
$ws.Range("I122").Value = 6311.4287
$ws.Range("J122").Value = 6833.3335
$ws.Range("K122").Value = 18934.2861
$ws.Range("L122").Value = 20500.0005
$ws.Range("M122").Value = -16484.2861
$ws.Range("N122").Value = -25400.0005

$ws.Range("H132").Value = 45999.8
$ws.Range("I132").Value = 44999.75
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 134999.25
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -132469.25
$ws.Range("N132").Value = -155060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 878.3
$ws.Range("I113").Value = 912
$ws.Range("J113").Value = 799.6667
$ws.Range("K113").Value = 2736
$ws.Range("L113").Value = 2399.0001
$ws.Range("M113").Value = -566
$ws.Range("N113").Value = -6739.0001

$ws.Range("H122").Value = 501724.5
$ws.Range("I122").Value = 667633
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 2002899
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -2000449
$ws.Range("N122").Value = -16897
